# Fixed naive component forecaster bug - Presentation state 11.02.
#
# A new (more recent) quarter of error data is inserted at the top of the
# table (row 2). The pre-existing rows shift down by one row, and the
# oldest row (which was row 11) falls off the bottom of the table, since
# the table keeps a fixed window of 10 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) values for columns B:G, rows 2-11, before
# overwriting anything, so the shift-down can be performed safely.
$oldValues = @{}
for ($r = 2; $r -le 11; $r++) {
    $oldValues[$r] = @(
        $ws.Cells.Item($r, 2).Value2,  # B
        $ws.Cells.Item($r, 3).Value2,  # C
        $ws.Cells.Item($r, 4).Value2,  # D
        $ws.Cells.Item($r, 5).Value2,  # E
        $ws.Cells.Item($r, 6).Value2,  # F
        $ws.Cells.Item($r, 7).Value2   # G
    )
}

# Shift rows 2-10 down into rows 3-11 (row 11's old data is discarded).
for ($r = 10; $r -ge 2; $r--) {
    $vals = $oldValues[$r]
    $dest = $r + 1
    $ws.Cells.Item($dest, 2).Value2 = $vals[0]
    $ws.Cells.Item($dest, 3).Value2 = $vals[1]
    $ws.Cells.Item($dest, 4).Value2 = $vals[2]
    $ws.Cells.Item($dest, 5).Value2 = $vals[3]
    $ws.Cells.Item($dest, 6).Value2 = $vals[4]
    $ws.Cells.Item($dest, 7).Value2 = $vals[5]
}

# Write the newly computed error metrics for the latest quarter into row 2.
$ws.Cells.Item(2, 2).Value2 = 0.2015370511150554
$ws.Cells.Item(2, 3).Value2 = 0.3501553535809984
$ws.Cells.Item(2, 4).Value2 = 0.2617601871928103
$ws.Cells.Item(2, 5).Value2 = 0.5116250455097076
$ws.Cells.Item(2, 6).Value2 = 0.4867637343656181
$ws.Cells.Item(2, 7).Value2 = 15
